$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 200
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = -188
$ws.Range("N6").Value = -1124
$ws.Range("H40").Value = 1360.591
$ws.Range("I40").Value = 1244.0834
$ws.Range("J40").Value = 1500.4
$ws.Range("K40").Value = 1244.0834
$ws.Range("L40").Value = 1500.4
$ws.Range("M40").Value = -1069.0834
$ws.Range("N40").Value = -1850.4
$ws.Range("H100").Value = 2764.5454
$ws.Range("I100").Value = 1481
$ws.Range("J100").Value = 3834.1667
$ws.Range("K100").Value = 1481
$ws.Range("L100").Value = 3834.1667
$ws.Range("M100").Value = -940
$ws.Range("N100").Value = -4916.1667
$ws.Range("H106").Value = 580
$ws.Range("I106").Value = 422.22223
$ws.Range("K106").Value = 422.22223
$ws.Range("M106").Value = 208.77777
$ws.Range("H116").Value = 16669450
$ws.Range("I116").Value = 40001800
$ws.Range("J116").Value = 3485.7144
$ws.Range("K116").Value = 40001800
$ws.Range("L116").Value = 3485.7144
$ws.Range("M116").Value = -39998358
$ws.Range("N116").Value = -10369.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2295.2222
$ws.Range("I61").Value = 1940.8
$ws.Range("J61").Value = 2738.25
$ws.Range("K61").Value = 1940.8
$ws.Range("L61").Value = 2738.25
$ws.Range("M61").Value = -1728.8
$ws.Range("N61").Value = -3162.25
$ws.Range("H63").Value = 3817.4443
$ws.Range("J63").Value = 2950
$ws.Range("L63").Value = 2950
$ws.Range("N63").Value = -4322
$ws.Range("H66").Value = 3817.4443
$ws.Range("J66").Value = 2950
$ws.Range("L66").Value = 14750
$ws.Range("N66").Value = -21614
$ws.Range("H74").Value = 1973.1428
$ws.Range("I74").Value = 2037.3334
$ws.Range("J74").Value = 1925
$ws.Range("K74").Value = 2037.3334
$ws.Range("L74").Value = 1925
$ws.Range("M74").Value = -1163.3334
$ws.Range("N74").Value = -3673
$ws.Range("H77").Value = 1973.1428
$ws.Range("I77").Value = 2037.3334
$ws.Range("J77").Value = 1925
$ws.Range("K77").Value = 10186.667
$ws.Range("L77").Value = 9625
$ws.Range("M77").Value = -5818.666999999999
$ws.Range("N77").Value = -18361
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 1366.1765
$ws.Range("I132").Value = 1072.2142
$ws.Range("J132").Value = 2738
$ws.Range("K132").Value = 3216.6426
$ws.Range("L132").Value = 8214
$ws.Range("M132").Value = -686.6425999999997
$ws.Range("N132").Value = -13274
$ws.Range("H136").Value = 2295.2222
$ws.Range("I136").Value = 1940.8
$ws.Range("J136").Value = 2738.25
$ws.Range("K136").Value = 5822.4
$ws.Range("L136").Value = 8214.75
$ws.Range("M136").Value = -3272.4
$ws.Range("N136").Value = -13314.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2611.0417
$ws.Range("I134").Value = 2220.842
$ws.Range("J134").Value = 4093.8
$ws.Range("K134").Value = 6662.526
$ws.Range("L134").Value = 12281.4
$ws.Range("M134").Value = -4127.526
$ws.Range("N134").Value = -17351.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 32857.145
$ws.Range("I3").Value = 50000
$ws.Range("J3").Value = 30000
$ws.Range("K3").Value = 50000
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = -49887
$ws.Range("N3").Value = -30226
$ws.Range("H16").Value = 1966.6666
$ws.Range("I16").Value = 1950
$ws.Range("K16").Value = 1950
$ws.Range("M16").Value = -1663
$ws.Range("H31").Value = 22225440
$ws.Range("I31").Value = 35716950
$ws.Range("K31").Value = 35716950
$ws.Range("M31").Value = -35716655
$ws.Range("H34").Value = 22225440
$ws.Range("I34").Value = 35716950
$ws.Range("K34").Value = 35716950
$ws.Range("M34").Value = -35716748
$ws.Range("H58").Value = 1633.5
$ws.Range("I58").Value = 1576.2273
$ws.Range("J58").Value = 1948.5
$ws.Range("K58").Value = 1576.2273
$ws.Range("L58").Value = 1948.5
$ws.Range("M58").Value = -1373.2273
$ws.Range("N58").Value = -2354.5
$ws.Range("H113").Value = 1966.6666
$ws.Range("I113").Value = 1950
$ws.Range("K113").Value = 1950
$ws.Range("M113").Value = 220
$ws.Range("H132").Value = 944.8889
$ws.Range("I132").Value = 958.25714
$ws.Range("K132").Value = 2874.77142
$ws.Range("M132").Value = -344.77142
$ws.Range("H134").Value = 949.2619
$ws.Range("I134").Value = 884.225
$ws.Range("J134").Value = 2250
$ws.Range("K134").Value = 2652.675
$ws.Range("L134").Value = 6750
$ws.Range("M134").Value = -117.6750000000002
$ws.Range("N134").Value = -11820
$ws.Range("H136").Value = 1633.5
$ws.Range("I136").Value = 1576.2273
$ws.Range("J136").Value = 1948.5
$ws.Range("K136").Value = 4728.6819
$ws.Range("L136").Value = 1948.5
$ws.Range("M136").Value = -2178.6819
$ws.Range("N136").Value = -10945.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 805201.6
$ws.Range("J12").Value = 966220.1
$ws.Range("L12").Value = 2898660.3
$ws.Range("N12").Value = -2899006.3
$ws.Range("H122").Value = 803.8919
$ws.Range("I122").Value = 497.94446
$ws.Range("J122").Value = 1093.7368
$ws.Range("K122").Value = 4481.50014
$ws.Range("L122").Value = 9843.6312
$ws.Range("M122").Value = -2031.50014
$ws.Range("N122").Value = -14743.6312
$ws.Range("H131").Value = 14495457
$ws.Range("I131").Value = 12853.75
$ws.Range("J131").Value = 16394815
$ws.Range("K131").Value = 38561.25
$ws.Range("L131").Value = 49184445
$ws.Range("M131").Value = -33521.25
$ws.Range("N131").Value = -49194525

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 13366.667
$ws.Range("J5").Value = 12800
$ws.Range("L5").Value = 12800
$ws.Range("N5").Value = -13024

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1729.6666
$ws.Range("I93").Value = 1644.5
$ws.Range("K93").Value = 1644.5
$ws.Range("M93").Value = -396.5
$ws.Range("H122").Value = 11908839
$ws.Range("I122").Value = 14709820
$ws.Range("J122").Value = 4666.5
$ws.Range("K122").Value = 44129460
$ws.Range("L122").Value = 13999.5
$ws.Range("M122").Value = -44127010
$ws.Range("N122").Value = -18899.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8931727
$ws.Range("I122").Value = 12502223
$ws.Range("K122").Value = 37506669
$ws.Range("M122").Value = -37504219
$ws.Range("H132").Value = 1226.5
$ws.Range("I132").Value = 1004.58826
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 3013.76478
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -483.76478
$ws.Range("N132").Value = -20057
